# Edit script: insert a newly-injured player (Rubtsov German, Spartak) into the
# "snapshot" sheet at its sorted position, refresh the "scraped_at" (column K)
# timestamps for every data row (the whole snapshot was re-scraped), and record
# the new injury on the "new_injured" sheet.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. "snapshot" sheet: insert the new row and refresh timestamps.
# ---------------------------------------------------------------------------
$snapshot = $wb.Worksheets.Item("snapshot")

# Make room for the new record right after the other "Spartak" (SPR) rows,
# before the first "Salavat Yulaev" (SYuL) row.
$snapshot.Rows.Item(34).Insert()

$newRow = @(
    "СПР",
    "Спартак",
    "spartak",
    "Рубцов Герман",
    "95",
    "нападающий",
    "22494",
    "1369_СПР_рубцовгерман",
    "injured_active",
    "https://www.khl.ru/clubs/spartak/team/",
    "2025-10-28T13:07:13.802963+00:00"
)
for ($col = 1; $col -le $newRow.Length; $col++) {
    $cell = $snapshot.Cells.Item(34, $col)
    # Force plain text so numeric-/date-looking values ("95", "22494", the
    # ISO timestamp) are not auto-converted to numbers/dates by Excel.
    $cell.Value = "'" + $newRow[$col - 1]
    $cell.ClearFormats()
}

# Refresh column K ("scraped_at") for every data row (2..50) to the values
# recorded by this re-scrape. Row 34 is the freshly-inserted row above; rows
# 2-33 keep their existing data but get a new scrape timestamp, and rows
# 35-50 are the former rows 34-49 shifted down by the insert.
$scrapedAt = @(
    "2025-10-28T13:06:33.727615+00:00",
    "2025-10-28T13:06:35.795830+00:00",
    "2025-10-28T13:06:35.795847+00:00",
    "2025-10-28T13:06:35.795855+00:00",
    "2025-10-28T13:06:35.795863+00:00",
    "2025-10-28T13:06:35.795870+00:00",
    "2025-10-28T13:06:38.040824+00:00",
    "2025-10-28T13:06:38.040840+00:00",
    "2025-10-28T13:06:40.166925+00:00",
    "2025-10-28T13:06:42.930049+00:00",
    "2025-10-28T13:06:42.930078+00:00",
    "2025-10-28T13:06:45.947213+00:00",
    "2025-10-28T13:06:45.947242+00:00",
    "2025-10-28T13:06:45.947260+00:00",
    "2025-10-28T13:06:45.947278+00:00",
    "2025-10-28T13:06:53.916643+00:00",
    "2025-10-28T13:06:56.705517+00:00",
    "2025-10-28T13:06:59.022961+00:00",
    "2025-10-28T13:07:01.338159+00:00",
    "2025-10-28T13:07:01.338189+00:00",
    "2025-10-28T13:07:03.607644+00:00",
    "2025-10-28T13:07:03.607678+00:00",
    "2025-10-28T13:07:03.607696+00:00",
    "2025-10-28T13:07:03.607714+00:00",
    "2025-10-28T13:07:03.607731+00:00",
    "2025-10-28T13:07:05.982882+00:00",
    "2025-10-28T13:07:10.984358+00:00",
    "2025-10-28T13:07:10.984388+00:00",
    "2025-10-28T13:07:10.984407+00:00",
    "2025-10-28T13:07:10.984423+00:00",
    "2025-10-28T13:07:13.802919+00:00",
    "2025-10-28T13:07:13.802946+00:00",
    "2025-10-28T13:07:13.802963+00:00",
    "2025-10-28T13:07:16.242499+00:00",
    "2025-10-28T13:07:16.242527+00:00",
    "2025-10-28T13:07:16.242544+00:00",
    "2025-10-28T13:07:16.242559+00:00",
    "2025-10-28T13:07:16.242577+00:00",
    "2025-10-28T13:07:16.242592+00:00",
    "2025-10-28T13:07:16.242607+00:00",
    "2025-10-28T13:07:16.242621+00:00",
    "2025-10-28T13:07:16.242636+00:00",
    "2025-10-28T13:07:18.675266+00:00",
    "2025-10-28T13:07:18.675293+00:00",
    "2025-10-28T13:07:23.522488+00:00",
    "2025-10-28T13:07:25.897226+00:00",
    "2025-10-28T13:07:25.897256+00:00",
    "2025-10-28T13:07:25.897274+00:00",
    "2025-10-28T13:07:25.897291+00:00"
)
for ($i = 0; $i -lt $scrapedAt.Length; $i++) {
    $row = $i + 2
    $cell = $snapshot.Cells.Item($row, 11)
    $cell.Value = "'" + $scrapedAt[$i]
    $cell.ClearFormats()
}

# ---------------------------------------------------------------------------
# 2. "new_injured" sheet: append the new injury record.
# ---------------------------------------------------------------------------
$newInjured = $wb.Worksheets.Item("new_injured")

$injuredRow = @(
    "СПР",
    "Спартак",
    "Рубцов Герман",
    "1369_СПР_рубцовгерман",
    "INJURED_NEW",
    "2025-10-28T21:07:26.410868+08:00",
    "2025-10-28"
)
for ($col = 1; $col -le $injuredRow.Length; $col++) {
    $cell = $newInjured.Cells.Item(2, $col)
    $cell.Value = "'" + $injuredRow[$col - 1]
    $cell.ClearFormats()
}

Write-Output "edit complete"
